# Update membership count values in column B per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 750
$ws.Range("B4").Value = 600
$ws.Range("B5").Value = 120
$ws.Range("B6").Value = 100
$ws.Range("B7").Value = 120
$ws.Range("B8").Value = 300
